$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.98
$ws.Range("H2").Value = 3.6
$ws.Range("I2").Value = 3.25
$ws.Range("J2").Value = 2.5
$ws.Range("L2").Value = 3.65
$ws.Range("M2").Value = 1.25
$ws.Range("N2").Value = 3.25
$ws.Range("O2").Value = 1.75
$ws.Range("P2").Value = 1.87
$ws.Range("Q2").Value = 2.72
$ws.Range("R2").Value = 1.35
$ws.Range("U2").Value = 1.65
$ws.Range("V2").Value = 1.98
$ws.Range("W2").Value = 8
$ws.Range("X2").Value = 9.75
$ws.Range("Y2").Value = 8.5
$ws.Range("Z2").Value = 17.5
$ws.Range("AA2").Value = 15
$ws.Range("AB2").Value = 25
$ws.Range("AC2").Value = 11.5
$ws.Range("AD2").Value = 7
$ws.Range("AE2").Value = 14
$ws.Range("AF2").Value = 60
$ws.Range("AG2").Value = 450
$ws.Range("AH2").Value = 10.75
$ws.Range("AI2").Value = 17.5
$ws.Range("AJ2").Value = 11.5
$ws.Range("AK2").Value = 45
$ws.Range("AL2").Value = 27
$ws.Range("AM2").Value = 35

# Row 5 updates
$ws.Range("O5").Value = 2.4
$ws.Range("P5").Value = 1.53
$ws.Range("AO5").Value = 8
$ws.Range("AP5").Value = 1.83
$ws.Range("AQ5").Value = 2.03
